{"js": "// Insert six new paragraphs right after the paragraph that ends with\n// \"Relation: Kind Statements aggregations. Transform: Kind Resources\n// related to themselves (ID), then Relations to other Resource via\n// Dataflow Kinds domain / range relationship (ordered).\"\n//\n// New paragraphs (in order):\n//   1. (blank)\n//   2. \"Relation: aggregated aligned entities. Views (transforms). Kind\n//       members occurring in Statement Resource(s). Functors / Monads:\"\n//   3. (blank)\n//   4. \"Relation<A>::flatMap(F : Function<A, B>) : Relation<B>;\"\n//   5. (blank)\n//   6. \"Function: declarative dataflow transform.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"Relation: Kind Statements aggregations. Transform: Kind Resources \" +\n  \"related to themselves (ID), then Relations to other Resource via \" +\n  \"Dataflow Kinds domain / range relationship (ordered).\";\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the anchor paragraph to insert after.\");\n}\n\nconst newTexts = [\n  \"\",\n  \"Relation: aggregated aligned entities. Views (transforms). Kind members occurring in Statement Resource(s). Functors / Monads:\",\n  \"\",\n  \"Relation<A>::flatMap(F : Function<A, B>) : Relation<B>;\",\n  \"\",\n  \"Function: declarative dataflow transform.\",\n];\n\n// Insert in order, each time right after the previously inserted paragraph\n// (or after the anchor for the first one), so the final order matches the\n// diff exactly.\nlet insertAfter = anchor;\nfor (const text of newTexts) {\n  insertAfter = insertAfter.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert six new paragraphs right after the paragraph that ends with\n# \"Relation: Kind Statements aggregations. Transform: Kind Resources\n# related to themselves (ID), then Relations to other Resource via\n# Dataflow Kinds domain / range relationship (ordered).\"\n#\n# New paragraphs (in order):\n#   1. (blank)\n#   2. \"Relation: aggregated aligned entities. Views (transforms). Kind\n#       members occurring in Statement Resource(s). Functors / Monads:\"\n#   3. (blank)\n#   4. \"Relation<A>::flatMap(F : Function<A, B>) : Relation<B>;\"\n#   5. (blank)\n#   6. \"Function: declarative dataflow transform.\"\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Relation: Kind Statements aggregations. Transform: Kind Resources related to themselves (ID), then Relations to other Resource via Dataflow Kinds domain / range relationship (ordered).\"\n\n# Confirm the anchor text is present (idiomatic Find-based check) before\n# doing the actual structural edit via the Paragraphs collection.\n$checkRange = $d.Content\n$found = $checkRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Could not find the anchor paragraph to insert after.\"\n}\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $ptext = $p.Range.Text.TrimEnd([char]13)\n    if ($ptext -eq $anchorText) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the anchor paragraph to insert after.\"\n}\n\n$newTexts = @(\n    \"\",\n    \"Relation: aggregated aligned entities. Views (transforms). Kind members occurring in Statement Resource(s). Functors / Monads:\",\n    \"\",\n    \"Relation<A>::flatMap(F : Function<A, B>) : Relation<B>;\",\n    \"\",\n    \"Function: declarative dataflow transform.\"\n)\n\n$cur = $target\nforeach ($t in $newTexts) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $cur.Next()\n    if ($t -ne \"\") {\n        $cur.Range.InsertAfter($t)\n    }\n}\n"}
